$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.034.75"
$ws.Range("E2").Value = "  +5.76%  "
$ws.Range("D3").Value = "3.719.00"
$ws.Range("E3").Value = "  +19.69%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.30"
$ws.Range("E5").Value = "  +7.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.50"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("D7").Value = "3.718.73"
$ws.Range("E7").Value = "  +19.73%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +5.56%  "
$ws.Range("E10").Value = "  +7.60%  "
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.503"
$ws.Range("E12").Value = "  +7.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.51"
$ws.Range("E13").Value = "  +12.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000255"
$ws.Range("E14").Value = "  +5.88%  "
$ws.Range("D15").Value = "4.337.06"
$ws.Range("E15").Value = "  +19.64%  "
$ws.Range("D16").Value = "3.715.99"
$ws.Range("E16").Value = "  +19.67%  "
$ws.Range("D17").Value = "71.105.72"
$ws.Range("E17").Value = "  +6.01%  "
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  +6.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "519.31"
$ws.Range("E20").Value = "  +5.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.91"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.33"
$ws.Range("E22").Value = "  +19.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.746"
$ws.Range("E23").Value = "  +8.66%  "
$ws.Range("E24").Value = "  +12.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.67"
$ws.Range("E25").Value = "  +6.06%  "
$ws.Range("E26").Value = "  +7.42%  "
$ws.Range("E27").Value = "  +10.97%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  +9.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.18"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  +11.29%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.02"
$ws.Range("E32").Value = "  +13.55%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000112"
$ws.Range("E33").Value = "  +18.59%  "
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +9.80%  "
$ws.Range("E37").Value = "  +9.81%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.347"
$ws.Range("E38").Value = "  +10.84%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  +11.60%  "
$ws.Range("E40").Value = "  +8.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.64"
$ws.Range("E41").Value = "  +5.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "434.52"
$ws.Range("E42").Value = "  +17.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.90"
$ws.Range("E43").Value = "  -5.46%  "
$ws.Range("D44").Value = "3.156.04"
$ws.Range("E44").Value = "  +12.86%  "
$ws.Range("E45").Value = "  +6.55%  "
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0367"
$ws.Range("E47").Value = "  +6.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.21"
$ws.Range("E48").Value = "  +10.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.81"
$ws.Range("E49").Value = "  +3.58%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.50"
$ws.Range("E51").Value = "  +9.69%  "
